$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rutas_registros")

# B2: keep the existing "Test" value but pick up the same direct formatting
# (style index) used elsewhere in column B. Clearing then re-writing makes the
# engine re-materialize the cell with column B's default formatting.
$ws.Range("B2").Clear()
$ws.Range("B2").Value = "Test"

# A3: was an empty numeric placeholder cell; now holds the text "20240810"
# (a date-like string, not a real number) while keeping its existing style.
# Assign it through a formula and then convert to a static value so the
# digit string isn't auto-coerced into a number.
$ws.Range("A3").Formula = '="20240810"'
$ws.Range("A3").Copy()
$ws.Range("A3").PasteSpecial(-4163)

# B3: brand new cell added next to A3, same "Test" label used in B2.
$ws.Range("B3").Value = "Test"
